$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 70
$ws.Range("I12").Value = 70
$ws.Range("K12").Value = 70
$ws.Range("M12").Value = 100

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1021.75
$ws.Range("I28").Value = 695.6667
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 695.6667
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -210.6667
$ws.Range("N28").Value = -2970

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1450
$ws.Range("I33").Value = 400.33334
$ws.Range("K33").Value = 400.33334
$ws.Range("M33").Value = -171.33334

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3650
$ws.Range("J64").Value = 3650
$ws.Range("L64").Value = 3650
$ws.Range("N64").Value = -4146

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3650
$ws.Range("J67").Value = 3650
$ws.Range("L67").Value = 3650
$ws.Range("N67").Value = -5366

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1715
$ws.Range("J70").Value = 1875
$ws.Range("L70").Value = 5625
$ws.Range("N70").Value = -6165

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1715
$ws.Range("J73").Value = 1875
$ws.Range("L73").Value = 5625
$ws.Range("N73").Value = -7497

# ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 43250
$ws.Range("J82").Value = 46500
$ws.Range("L82").Value = 139500
$ws.Range("N82").Value = -140312

# ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 43250
$ws.Range("J85").Value = 46500
$ws.Range("L85").Value = 139500
$ws.Range("N85").Value = -142308

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 161
$ws.Range("I107").Value = 161
$ws.Range("K107").Value = 161
$ws.Range("M107").Value = 1759

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1999.9
$ws.Range("I113").Value = 1999.875
$ws.Range("K113").Value = 1999.875
$ws.Range("M113").Value = 1254.125

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9296.25
$ws.Range("I116").Value = 15152.5
$ws.Range("J116").Value = 3440
$ws.Range("K116").Value = 15152.5
$ws.Range("L116").Value = 3440
$ws.Range("M116").Value = -11710.5
$ws.Range("N116").Value = -10324

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4687.222
$ws.Range("I132").Value = 4687.222
$ws.Range("K132").Value = 14061.666
$ws.Range("M132").Value = -11531.666

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -3754

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 705.5
$ws.Range("I110").Value = 705.5
$ws.Range("K110").Value = 705.5
$ws.Range("M110").Value = 1339.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 680.2
$ws.Range("I94").Value = 367
$ws.Range("K94").Value = 367
$ws.Range("M94").Value = 84

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15072.5
$ws.Range("I105").Value = 16097.272
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 16097.272
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -14350.272
$ws.Range("N105").Value = -7294

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3858.625
$ws.Range("J107").Value = 1875
$ws.Range("L107").Value = 1875
$ws.Range("N107").Value = -5715

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1834.3334
$ws.Range("J94").Value = 1036.6666
$ws.Range("L94").Value = 1036.6666
$ws.Range("N94").Value = -1938.6666

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 974.2308
$ws.Range("I105").Value = 1136.1111
$ws.Range("J105").Value = 610
$ws.Range("K105").Value = 1136.1111
$ws.Range("L105").Value = 610
$ws.Range("M105").Value = 610.8888999999999
$ws.Range("N105").Value = -4104

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6568
$ws.Range("I134").Value = 1556
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 4668
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -2133
$ws.Range("N134").Value = -29070

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1582.75
$ws.Range("J131").Value = 1966.6
$ws.Range("L131").Value = 5899.799999999999
$ws.Range("N131").Value = -15979.8

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 450
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 450
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4290

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5783.778
$ws.Range("I113").Value = 5295.5713
$ws.Range("J113").Value = 7492.5
$ws.Range("K113").Value = 5295.5713
$ws.Range("L113").Value = 7492.5
$ws.Range("M113").Value = -3125.5713
$ws.Range("N113").Value = -11832.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 475
$ws.Range("I22").Value = 475
$ws.Range("K22").Value = 475
$ws.Range("M22").Value = -180

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 475
$ws.Range("I27").Value = 475
$ws.Range("K27").Value = 475
$ws.Range("M27").Value = -368

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5500
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -7498

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5500
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 30000
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -37488

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4837
$ws.Range("I126").Value = 3116
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 9348
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -6878
$ws.Range("N126").Value = -34940
